# Migração de OpenShift para Firebase
# Update the big "App / Clube do Livro" logo text box on slide 1:
#  - reflow the text into three separate centered paragraphs
#  - bump the font size way up and make it bold
#  - reposition/resize the text box to its new (larger) bounding box

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(3)

# Reflow text: "App " / "Clube do " / "Livro" on three separate lines.
$shp.TextFrame.TextRange.Text = "App `rClube do `rLivro"

# New formatting: 150pt bold text (was 80pt, not bold).
$shp.TextFrame.TextRange.Font.Size = 150
$shp.TextFrame.TextRange.Font.Bold = $true

# New position/size (EMU -> points, 12700 EMU per point).
$shp.Left = 804233 / 12700
$shp.Top = -99392 / 12700
$shp.Width = 7643439 / 12700
$shp.Height = 7017306 / 12700
